# ArchitectureDiagram.pptx update
#   - "Model" rectangle label becomes "Model(s)"
#   - the "Web" cloud shape and the elbow connector feeding it are removed
#
# (Per-slide/layout/master "datetimeFigureOut" footer fields are an
# auto-updating PowerPoint field that the host recomputes on save; they
# are not touched here.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $targetId) {
  for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Id -eq $targetId) {
      return $sh
    }
  }
  return $null
}

# 1) "Model" -> "Model(s)" (shape id 5, "Rectangle 45")
$modelShape = Get-ShapeById $s 5
$modelShape.TextFrame.TextRange.Text = "Model(s)"

# 2) Remove the "Web" cloud (id 51) and its elbow connector (id 52).
#    Delete the connector first since it references the cloud shape.
$connector = Get-ShapeById $s 52
if ($connector -ne $null) {
  $connector.Delete()
}

$cloud = Get-ShapeById $s 51
if ($cloud -ne $null) {
  $cloud.Delete()
}
